$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared strings must be created in this order so the new <si> entries
# --- land at the same indices the target workbook uses (12..16).
$ws.Range("C9").Value  = "Gemaakte klassendiagram implementeren"   # -> shared string 12
$ws.Range("C8").Value  = "PLSQL SP's maken"                         # -> shared string 13
$ws.Range("C7").Value  = "PLSQL implementeren in ASP.NET"          # -> shared string 14
$ws.Range("C10").Value = "Code  implementeren"                     # -> shared string 15
$ws.Range("C16").Value = "Testen"                                  # -> shared string 16

# --- Row 7: date moved from 28 Jun 2015 to 30 Jun 2015
$ws.Range("A7").Value = 42185

# --- Rows 11,12,14: fill in the member names (B column) for the days
# --- that already existed as blank rows.
$ws.Range("B11").Value = "Kees W."
$ws.Range("B12").Value = "Stan W."
$ws.Range("B14").Value = "Kees W."

# --- Row 10: new day block (1 Jul 2015) - give A10 the same date format as A7/A4
$ws.Range("A4").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 42186
$ws.Range("B10").Value = "Mark C."

# --- Row 13: new day block (2 Jul 2015) - give A13 the same date format as A7/A4
$ws.Range("A4").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 42187
$ws.Range("B13").Value = "Mark C."

# --- Fill in column C ("Uit te voeren activiteit") for rows 11-15 (all "Code  implementeren")
$ws.Range("C11").Value = "Code  implementeren"
$ws.Range("C12").Value = "Code  implementeren"
$ws.Range("C13").Value = "Code  implementeren"
$ws.Range("C14").Value = "Code  implementeren"

# --- Row 15 (new): Stan W., Code implementeren - reuse formatting from row 14 (A/B cols)
$ws.Range("A14:B14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("B15").Value = "Stan W."
$ws.Range("C15").Value = "Code  implementeren"

# --- Row 16 (new): 3 Jul 2015, Mark C., Testen
$ws.Range("A13:B13").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 42188
$ws.Range("B16").Value = "Mark C."
$ws.Range("C16").Value = "Testen"

# --- Row 17 (new): Kees W., Testen
$ws.Range("A14:B14").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B17").Value = "Kees W."
$ws.Range("C17").Value = "Testen"

# --- Row 18 (new): Stan W., Testen - no A18 cell on this row
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Stan W."
$ws.Range("C18").Value = "Testen"

# --- Row 19 (new): closing rule across the table (thin top border)
$ws.Range("A19:G19").Borders.Item(8).LineStyle = 1

$ws.Range("A1:G19").Select()
